# Auto-generated edit script
# Applies the updated Step1_Data raw values and their downstream
# recomputed cumulative-sum (Step2_Sj) and threshold crossing (Step3_DataPts_*)
# results, per the target diff.

$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$updates = @{
    "D2" = 0.1394125610291484
    "F2" = 0.3513066017617835
    "G2" = 0.1254598515707391
    "J2" = 0.0300724787281603
    "L2" = 0.01847361026986566
    "M2" = 0.007290732174072548
    "O2" = 0.1345315493922528
    "P2" = 0.05155416210444903
    "R2" = 0.01016696765940623
    "S2" = 0.006832566738613331
    "T2" = 0.03051147419784146
    "U2" = 0.06775426251041856
    "AF2" = 0.0052268959931728
    "AI2" = 0.02140628587007671
    "D3" = 0.2230528271171323
    "E3" = 0.01467447455852309
    "F3" = 0.3180285973369705
    "G3" = 0.1234638801958112
    "H3" = 0.007870437627838098
    "I3" = 0.003290327977106782
    "O3" = 0.07536738752198818
    "P3" = 0.01858482704631432
    "T3" = 0.05536532651630817
    "U3" = 0.08156858624900915
    "AC3" = 0.005035748639735089
    "AH3" = 0.004966203497648527
    "AI3" = 0.06873137571561452
    "D4" = 0.01793053357807211
    "E4" = 0.1913214668559838
    "F4" = 0.1961249494970221
    "G4" = 0.3657383938457505
    "I4" = 0.05314496981627796
    "O4" = 0.01881655585632886
    "P4" = 0.03449847998565762
    "S4" = 0.006895444680382039
    "T4" = 0.01213941746729117
    "U4" = 0.05349163077447195
    "AC4" = 0.0108621898849717
    "AG4" = 0.006329992897913441
    "AI4" = 0.03270597485987679
    "D5" = 0.0166222887523202
    "E5" = 0.2317556573435376
    "F5" = 0.1843764365560175
    "G5" = 0.3793177551174527
    "I5" = 0.03959424036719365
    "K5" = 0.0007034522308447797
    "O5" = 0.001027404746738711
    "P5" = 0.01175846956979458
    "T5" = 0.01924126858972292
    "U5" = 0.06063891458070322
    "AI5" = 0.05496411214567382
    "D6" = 0.1319855251659318
    "E6" = 0.04579184825827361
    "F6" = 0.3347741234607933
    "G6" = 0.2691990305626665
    "H6" = 0.002868226877661598
    "I6" = 0.02710799338392731
    "J6" = 0.005688870644421639
    "O6" = 0.02106502934253471
    "P6" = 0.006824756269014002
    "T6" = 0.03999098481229016
    "U6" = 0.04462037041324453
    "AA6" = 0.008009720156683404
    "AC6" = 0.008673289592570717
    "AI6" = 0.05340023105998669
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$updates = @{
    "D2" = 0.1394125610291484
    "E2" = 0.1394125610291484
    "F2" = 0.4907191627909319
    "G2" = 0.616179014361671
    "H2" = 0.616179014361671
    "I2" = 0.616179014361671
    "J2" = 0.6462514930898313
    "K2" = 0.6462514930898313
    "L2" = 0.6647251033596969
    "M2" = 0.6720158355337694
    "N2" = 0.6720158355337694
    "O2" = 0.8065473849260222
    "P2" = 0.8581015470304713
    "Q2" = 0.8581015470304713
    "R2" = 0.8682685146898775
    "S2" = 0.8751010814284909
    "T2" = 0.9056125556263324
    "U2" = 0.973366818136751
    "V2" = 0.973366818136751
    "W2" = 0.973366818136751
    "X2" = 0.973366818136751
    "Y2" = 0.973366818136751
    "Z2" = 0.973366818136751
    "AA2" = 0.973366818136751
    "AB2" = 0.973366818136751
    "AC2" = 0.973366818136751
    "AD2" = 0.973366818136751
    "AE2" = 0.973366818136751
    "AF2" = 0.9785937141299238
    "AG2" = 0.9785937141299238
    "AH2" = 0.9785937141299238
    "D3" = 0.2230528271171323
    "E3" = 0.2377273016756554
    "F3" = 0.5557558990126259
    "G3" = 0.6792197792084371
    "H3" = 0.6870902168362752
    "I3" = 0.6903805448133821
    "J3" = 0.6903805448133821
    "K3" = 0.6903805448133821
    "L3" = 0.6903805448133821
    "M3" = 0.6903805448133821
    "N3" = 0.6903805448133821
    "O3" = 0.7657479323353702
    "P3" = 0.7843327593816846
    "Q3" = 0.7843327593816846
    "R3" = 0.7843327593816846
    "S3" = 0.7843327593816846
    "T3" = 0.8396980858979928
    "U3" = 0.921266672147002
    "V3" = 0.921266672147002
    "W3" = 0.921266672147002
    "X3" = 0.921266672147002
    "Y3" = 0.921266672147002
    "Z3" = 0.921266672147002
    "AA3" = 0.921266672147002
    "AB3" = 0.921266672147002
    "AC3" = 0.9263024207867371
    "AD3" = 0.9263024207867371
    "AE3" = 0.9263024207867371
    "AF3" = 0.9263024207867371
    "AG3" = 0.9263024207867371
    "AH3" = 0.9312686242843856
    "D4" = 0.01793053357807211
    "E4" = 0.2092520004340559
    "F4" = 0.405376949931078
    "G4" = 0.7711153437768284
    "H4" = 0.7711153437768284
    "I4" = 0.8242603135931064
    "J4" = 0.8242603135931064
    "K4" = 0.8242603135931064
    "L4" = 0.8242603135931064
    "M4" = 0.8242603135931064
    "N4" = 0.8242603135931064
    "O4" = 0.8430768694494353
    "P4" = 0.8775753494350929
    "Q4" = 0.8775753494350929
    "R4" = 0.8775753494350929
    "S4" = 0.8844707941154749
    "T4" = 0.8966102115827661
    "U4" = 0.950101842357238
    "V4" = 0.950101842357238
    "W4" = 0.950101842357238
    "X4" = 0.950101842357238
    "Y4" = 0.950101842357238
    "Z4" = 0.950101842357238
    "AA4" = 0.950101842357238
    "AB4" = 0.950101842357238
    "AC4" = 0.9609640322422097
    "AD4" = 0.9609640322422097
    "AE4" = 0.9609640322422097
    "AF4" = 0.9609640322422097
    "AG4" = 0.9672940251401232
    "AH4" = 0.9672940251401232
    "AI4" = 1
    "D5" = 0.0166222887523202
    "E5" = 0.2483779460958578
    "F5" = 0.4327543826518753
    "G5" = 0.8120721377693281
    "H5" = 0.8120721377693281
    "I5" = 0.8516663781365217
    "J5" = 0.8516663781365217
    "K5" = 0.8523698303673665
    "L5" = 0.8523698303673665
    "M5" = 0.8523698303673665
    "N5" = 0.8523698303673665
    "O5" = 0.8533972351141051
    "P5" = 0.8651557046838997
    "Q5" = 0.8651557046838997
    "R5" = 0.8651557046838997
    "S5" = 0.8651557046838997
    "T5" = 0.8843969732736227
    "U5" = 0.9450358878543259
    "V5" = 0.9450358878543259
    "W5" = 0.9450358878543259
    "X5" = 0.9450358878543259
    "Y5" = 0.9450358878543259
    "Z5" = 0.9450358878543259
    "AA5" = 0.9450358878543259
    "AB5" = 0.9450358878543259
    "AC5" = 0.9450358878543259
    "AD5" = 0.9450358878543259
    "AE5" = 0.9450358878543259
    "AF5" = 0.9450358878543259
    "AG5" = 0.9450358878543259
    "AH5" = 0.9450358878543259
    "AI5" = 0.9999999999999998
    "D6" = 0.1319855251659318
    "E6" = 0.1777773734242054
    "F6" = 0.5125514968849987
    "G6" = 0.7817505274476652
    "H6" = 0.7846187543253268
    "I6" = 0.8117267477092541
    "J6" = 0.8174156183536757
    "K6" = 0.8174156183536757
    "L6" = 0.8174156183536757
    "M6" = 0.8174156183536757
    "N6" = 0.8174156183536757
    "O6" = 0.8384806476962104
    "P6" = 0.8453054039652245
    "Q6" = 0.8453054039652245
    "R6" = 0.8453054039652245
    "S6" = 0.8453054039652245
    "T6" = 0.8852963887775146
    "U6" = 0.9299167591907591
    "V6" = 0.9299167591907591
    "W6" = 0.9299167591907591
    "X6" = 0.9299167591907591
    "Y6" = 0.9299167591907591
    "Z6" = 0.9299167591907591
    "AA6" = 0.9379264793474426
    "AB6" = 0.9379264793474426
    "AC6" = 0.9465997689400133
    "AD6" = 0.9465997689400133
    "AE6" = 0.9465997689400133
    "AF6" = 0.9465997689400133
    "AG6" = 0.9465997689400133
    "AH6" = 0.9465997689400133
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$updates = @{
    "D2" = 6
    "F2" = 0.616179014361671
    "G2" = 5
    "F3" = 0.5557558990126259
    "F4" = 0.7711153437768284
    "F5" = 0.8120721377693281
    "F6" = 0.5125514968849987
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$updates = @{
    "D2" = 14
    "F2" = 0.8065473849260222
    "G2" = 13
    "D3" = 14
    "F3" = 0.7657479323353702
    "G3" = 13
    "F4" = 0.7711153437768284
    "F5" = 0.8120721377693281
    "F6" = 0.7817505274476652
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$updates = @{
    "F2" = 0.8065473849260222
    "D3" = 19
    "F3" = 0.8396980858979928
    "G3" = 18
    "D4" = 8
    "F4" = 0.8242603135931064
    "G4" = 6
    "F5" = 0.8120721377693281
    "D6" = 8
    "F6" = 0.8117267477092541
    "G6" = 7
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$updates = @{
    "D2" = 19
    "F2" = 0.9056125556263324
    "G2" = 18
    "F3" = 0.921266672147002
    "D4" = 20
    "F4" = 0.950101842357238
    "G4" = 18
    "D5" = 20
    "F5" = 0.9450358878543259
    "G5" = 18
    "D6" = 20
    "F6" = 0.9299167591907591
    "G6" = 19
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

